$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: add "How to run" text for the GCN (first) row
$ws.Range("N11").Value = "Folder: TGCN/T-GCN/T-GCN-PyTorch:  python main.py --data shenzhen --model_name GCN --max_epochs 100 --learning_rate 0.0001 --weight_decay 0 --batch_size 32 --hidden_dim 64 --loss mse_with_regularizer --settings supervised --gpus 0 --seq_len 32 --pre_len 7"

# Row 13: GCN
$ws.Range("A13").Value = "GCN"
$ws.Range("B13").Value = 32
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 0.01
$ws.Range("E13").Value = 64
$ws.Range("F13").Value = 2016
$ws.Range("G13").Value = 23516632
$ws.Range("H13").Value = 4849
$ws.Range("I13").Value = 0.71
$ws.Range("J13").Value = 0.516
$ws.Range("K13").Value = 23516632
$ws.Range("L13").Value = 3000

# Row 14: STGCN
$ws.Range("A14").Value = "STGCN"
$ws.Range("B14").Value = 32
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = 0.001
$ws.Range("F14").Value = 707
$ws.Range("G14").Formula = "=1969.6^2"
$ws.Range("H14").Value = 1969.65
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = "python main.py --enable_cuda True --dataset pemsd7-m --n_his 32 --n_pred 7"

# Row 15: STGCN
$ws.Range("A15").Value = "STGCN"
$ws.Range("B15").Value = 32
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 0.01
$ws.Range("F15").Value = 769.7
$ws.Range("G15").Formula = "=2086.3^2"
$ws.Range("H15").Value = 2086.3
$ws.Range("L15").Value = 1000

# Row 16: STGCN
$ws.Range("A16").Value = "STGCN"
$ws.Range("B16").Value = 32
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 0.0001
$ws.Range("F16").Value = 625.493726
$ws.Range("G16").Formula = "=1831.15^2"
$ws.Range("H16").Value = 1831.15
$ws.Range("L16").Value = 1000

# Row 17: STGCN
$ws.Range("A17").Value = "STGCN"
$ws.Range("B17").Value = 32
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 0.0005
$ws.Range("F17").Value = 845.759016
$ws.Range("G17").Formula = "=2185^2"
$ws.Range("H17").Value = 2185
$ws.Range("L17").Value = 1000

# Row 18: STALSTM
$ws.Range("A18").Value = "STALSTM"
$ws.Range("B18").Value = 32
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 0.05
$ws.Range("F18").Value = 3356
$ws.Range("G18").Value = 79586797
$ws.Range("H18").Value = 4996
$ws.Range("L18").Value = 100
$ws.Range("N18").Value = "python main.py"

# Update selection to match target
$ws.Range("N23").Select() | Out-Null
